$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora/hour): every data row 2-51 advances from 15 to 16
$ws.Range("G2:G51").Value = "'16"
$ws.Range("G2:G51").Style = "Normal"

# Column D (Price) and E (Volume 1h %): refreshed quote snapshot per coin
# Row 2
$ws.Range("D2").Value = "'275.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.41%"
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = "'27.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.86%"
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").Value = "'4.762"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.53%"
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = "'0.06335"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.81%"
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = "'6.942"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.42%"
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("D7").Value = "'1.444"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'50.96%"
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").Value = "'0.8764"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.87%"
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").Value = "'0.1522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.80%"
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").Value = "'0.05049"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.96%"
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").Value = "'0.07495"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.09%"
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").Value = "'0.02871"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-7.47%"
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("D13").Value = "'0.09037"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.36%"
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").Value = "'0.001567"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.19%"
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").Value = "'0.0006358"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.75%"
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").Value = "'0.005779"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.76%"
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").Value = "'3.448"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.04%"
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").Value = "'3.301"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.34%"
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").Value = "'2.272"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.17%"
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("E20").Value = "'0.80%"
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("E21").Value = "'2.60%"
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = "'3.917"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.48%"
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").Value = "'0.04390"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.34%"
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").Value = "'0.001171"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.29%"
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").Value = "'0.003838"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.54%"
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("E26").Value = "'0.17%"
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("D27").Value = "'0.0001936"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'14.48%"
$ws.Range("E27").Style = "Normal"
# Row 40
$ws.Range("D40").Value = "'0.04117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.74%"
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = "'0.006854"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.28%"
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").Value = "'0.1174"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.33%"
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").Value = "'0.002050"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-12.99%"
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = "'0.01147"
$ws.Range("D44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = "'0.00005175"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.31%"
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("E46").Value = "'-36.83%"
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("D47").Value = "'0.02000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-11.36%"
$ws.Range("E47").Style = "Normal"
